$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up malformed French-locale numeric strings (non-breaking-space decimal
# separators like "3 ,5", stray typos like "0,3à" / "0,2§") that had been
# imported as text (shared strings) instead of numbers. Replace each with the
# correct numeric value so the column is fully numeric for the regression /
# t-test analysis.
$ws.Range("D38").Value  = 3.5
$ws.Range("B44").Value  = 0.28
$ws.Range("B82").Value  = 0.24
$ws.Range("B102").Value = 0.32
$ws.Range("D104").Value = 3.6
$ws.Range("D114").Value = 3.1
$ws.Range("B122").Value = 0.3
$ws.Range("B128").Value = 0.25
$ws.Range("D146").Value = 2.9
$ws.Range("B165").Value = 0.3

# Fix a rounding artifact in an existing numeric weight value.
$ws.Range("B158").Value = 0.28

# Update the view's scroll position / selection to match where the author
# left off after scrolling down to review the newly-cleaned rows.
$excel.ActiveWindow.ScrollRow = 152
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B179").Select()
